# Generate Report for Handoff
# Advances the "b.md" row from "Handed back: in sync with en-US" to
# "Ready for handoff" on all three sheets, recording the newly generated
# handoff xlf files / timestamps and the new error detail message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd445b9cf706317958a820f4caf95a8e7ecacb90/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/238856dc4c5767347abe10ee9cce28dab1dacd32/e2e/b.md."

# ---- Overview sheet : row 3 is b.md ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-10-27 06:28:09"

# ---- zh-cn sheet : row 3 is b.md ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces a genuine text value instead of auto-converting
# to a native boolean; resetting the style afterwards drops the resulting
# quote-prefix formatting so the cell style matches the original ("s=0").
$zh.Range("F3").Value = "'False"
$zh.Range("F3").Style = "Normal"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-10-27 06:27:56"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 40

# ---- de-de sheet : row 3 is b.md ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "'False"
$de.Range("F3").Style = "Normal"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-10-27 06:28:09"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 40
